# "fixing ACC for slide, calendar and dropdown"
#
# The three username cells (A3, A5, A8) hold the tester's email address
# "aya@4mail.top". The fix re-enters it with the proper capitalization
# ("Aya@4mail.top") and turns each occurrence into a mailto: hyperlink
# (Excel auto-applies its built-in "Hyperlink" cell style/font when this
# happens), then leaves the selection on A8 where the work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newEmail = "Aya@4mail.top"
$mailTarget = "mailto:" + $newEmail

$cells = @("A3", "A5", "A8")

foreach ($addr in $cells) {
    $rng = $ws.Range($addr)
    # Clear out the old shared-string text first so the stale
    # lower-case string isn't left dangling unused in the table,
    # then write back the corrected capitalization.
    $rng.ClearContents()
    $rng.Value = $newEmail
    $ws.Hyperlinks.Add($rng, $mailTarget)
}

# Leave the selection where editing finished.
[void]$ws.Range("A8").Select()

Write-Output "done"
